# "Added ChatApp MAUI data"
# - Adds the .NET MAUI package-size data point (30.7) to the "Package Size"
#   worksheet table (cell C4).
# - Recolors the bar-chart series fill from blue (0070C0) to purple (7030A0).
# - Updates the sheet's active selection to H23 (the last thing the author
#   clicked on before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Package Size")

# Add the missing "Size (MBs)" value for the .NET MAUI row.
$ws.Range("C4").Value = 30.7

# Recolor the column chart's single series (was solid blue 0070C0).
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Format.Fill.ForeColor.RGB = 10498160  # RGB(0x70, 0x30, 0xA0) = 7030A0

# Leave the selection where the author ended up.
$ws.Range("H23").Select()
